$d = $word.ActiveDocument

$values = @("Copper", "Nickel", "Chloramphenicol", "Ampicillin", "Metaldehyde", "Atrazine", "Tebuconazole", "Azoxystrobin")

$table = $d.Tables.Item(1)

for ($i = 1; $i -lt $values.Length; $i++) {
    $cell = $table.Cell($i, 2)
    $cell.Range.Text = $values[$i - 1]
}

# Last row: the real Word "_GoBack" bookmark tracks the most recent edit
# location, so typing into the final cell relocates it there (from the
# empty paragraph that used to follow the table) instead of leaving it in
# place. Reproduce that by moving the bookmark into the last cell, right
# after the newly typed text.
$lastCell = $table.Cell(8, 2)
$lastValue = $values[7]

# Type the text plus a throwaway sentinel character so the insertion point
# we bookmark is not the very last character in the paragraph (a boundary
# case the host mishandles for zero-length bookmarks).
$lastCell.Range.InsertAfter($lastValue + "X")

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$insertPos = $lastCell.Range.Start + $lastValue.Length
$target = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $target)

# Remove the sentinel now that the bookmark is safely anchored just before it.
$sentinel = $d.Range($insertPos, $insertPos + 1)
$sentinel.Delete()
